# Append a new order line (row 9) to the Regional Distributors order sheet.
# Existing rows store every value as text (inlineStr), even numeric-looking
# ones like SKU/Quantity/Cost, so force the new cells to Text format before
# assigning values to avoid Excel auto-coercing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A9:E9")
$newRow.NumberFormat = "@"

$ws.Range("A9").Value = "33576"
$ws.Range("B9").Value = "Cup - Cold (12oz)"
$ws.Range("C9").Value = "1"
$ws.Range("D9").Value = "69.75"
$ws.Range("E9").Value = "69.75"
